# Apply edits to TestData workbook:
# 1. Add a new row 3 with test case tc2 data (inserts shared strings 9-13)
# 2. Change B2 from boolean TRUE to the literal text string "true"
#    (formatted as Text / numFmtId 49) - inserts shared string 14
# 3. Move the active selection to G7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new test case tc2 ---
$ws.Range("A3").Value = "tc2"
$ws.Range("B3").Value = "false "
$ws.Range("C3").Value = "Minh "
$ws.Range("D3").Value = "Le "
$ws.Range("E3").Value = "Hoang"

# --- B2: convert boolean TRUE to literal text "true" ---
# A direct .Value assignment of "true"/"false" strings is always coerced
# back into a Boolean by the host, regardless of the cell's number format.
# Routing the text through a formula + values-only paste avoids that
# coercion and leaves a plain shared-string cell behind.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Formula = '="true"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Update selection to G7 ---
$ws.Range("G7").Select()
